# Apply the data update described by the commit:
# "updating model and datasets (auto generated commit)"
#
# This adjusts the Indiana COVID deaths-by-date-by-age-group sheet:
#  - a handful of existing daily covid_deaths counts are revised upward
#    (retroactive data corrections), and
#  - four new (date, agegrp, covid_deaths) rows are inserted, extending
#    the dataset (including a brand-new date, 2021-03-18 / serial 44275,
#    and a new trailing row for 2021-03-22 / serial 44279, 80+ age group).
#
# Work from the bottom of the sheet upward so row-number anchors used
# below (taken from the ORIGINAL/pre-edit layout) stay valid as each
# insertion shifts only the rows below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "YYYY-MM-DD HH:MM:SS"

function Insert-DataRow($rowIndex, $dateSerial, $ageGroup, $deaths) {
    $ws.Rows($rowIndex).Insert()
    $ws.Cells.Item($rowIndex, 1).Value = $dateSerial
    $ws.Cells.Item($rowIndex, 1).NumberFormat = $dateFmt
    $ws.Cells.Item($rowIndex, 2).Value = $ageGroup
    $ws.Cells.Item($rowIndex, 3).Value = $deaths
}

# --- New trailing row: 2021-03-22 (44279), 80+, 2 deaths ---
Insert-DataRow 1556 44279 "80+" 2

# --- New row: 2021-03-21 (44278), 60-69, 1 death ---
Insert-DataRow 1553 44278 "60-69" 1

# --- Revise 2021-03-19 (44276), 70-79: 2 -> 3 ---
$ws.Cells.Item(1547, 3).Value = 3

# --- New row: 2021-03-18 (44275), 60-69, 1 death (previously no data for this date) ---
Insert-DataRow 1546 44275 "60-69" 1

# --- Revise 2021-03-16 (44273), 80+: 7 -> 8 ---
$ws.Cells.Item(1544, 3).Value = 8

# --- Revise 2021-03-15 (44272), 80+: 3 -> 4 ---
$ws.Cells.Item(1539, 3).Value = 4

# --- Revise 2021-03-11 (44268), 70-79: 1 -> 2 ---
$ws.Cells.Item(1525, 3).Value = 2

# --- Revise 2021-03-10 (44267), 80+: 5 -> 6 ---
$ws.Cells.Item(1522, 3).Value = 6

# --- Revise 2021-03-10 (44267), 50-59: 3 -> 4 ---
$ws.Cells.Item(1519, 3).Value = 4

# --- Revise 2021-03-09 (44266), 80+: 2 -> 3 ---
$ws.Cells.Item(1518, 3).Value = 3

# --- Revise 2021-03-09 (44266), 70-79: 2 -> 3 ---
$ws.Cells.Item(1517, 3).Value = 3

# --- Revise 2021-03-08 (44265), 80+: 6 -> 7 ---
$ws.Cells.Item(1514, 3).Value = 7

# --- Revise 2021-03-05 (44262), 80+: 5 -> 6 ---
$ws.Cells.Item(1501, 3).Value = 6

# --- New row: 2021-02-28 (44257), 40-49, 1 death ---
Insert-DataRow 1478 44257 "40-49" 1

# --- Revise 2021-02-25 (44254), 70-79: 4 -> 5 ---
$ws.Cells.Item(1469, 3).Value = 5

# --- Revise 2021-02-12 (44241), 80+: 9 -> 10 ---
$ws.Cells.Item(1418, 3).Value = 10

# --- Revise 2021-02-12 (44241), 60-69: 3 -> 4 ---
$ws.Cells.Item(1416, 3).Value = 4

# --- Revise 2021-01-30 (44228), 40-49: 1 -> 2 ---
$ws.Cells.Item(1356, 3).Value = 2

# --- Revise 2021-01-17 (44215), 60-69: 10 -> 11 ---
$ws.Cells.Item(1294, 3).Value = 11
